# results_crossentropy.xlsx update:
#  - rename shared string "Type - Logical Augmentation VILT new way"
#    -> "Type - Logical Extension VILT new way" (row 81 label, via text write)
#  - clear the old, never-filled-in "VSR Logically Extended + Augmentation VILT"
#    block that lived at rows 71-79 (all zeros / #DIV/0! placeholders)
#  - re-add that block (now with real recorded results) plus the
#    "Contrastive VILT new way 0.01" block, at rows 86-94 (after the
#    existing "Logical Augmentation/Extension VILT new way" block at 81-84)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3 trials both models")

# --- rename the row-81 block label (shared string used elsewhere too) ---
$ws.Range("A81").Value = "Type - Logical Extension VILT new way"

# --- clear out the stale placeholder block that used to sit at rows 71-79 ---
$ws.Range("A71:I79").ClearContents()

# --- rows 86-94: header block + two filled-in result blocks ---

# Row 86: header row for "VSR Logically Extended + Augmentation VILT"
$ws.Range("A86").Value = "Type - VSR Logically Extended + Augmentation VILT"
$ws.Range("B86").Value = "Test acc."
$ws.Range("C86").Value = "Dev acc."
$ws.Range("D86").Value = "Min test"
$ws.Range("E86").Value = "Max test"
$ws.Range("F86").Value = "Mean test"
$ws.Range("G86").Value = "Min dev"
$ws.Range("H86").Value = "Max dev"
$ws.Range("I86").Value = "Mean dev"

# Row 87: Run#1
$ws.Range("A87").Value = "Run#1"
$ws.Range("B87").Value = 74.05
$ws.Range("C87").Value = 75.57
$ws.Range("D87").Formula = "=MIN(B87:B91)"
$ws.Range("E87").Formula = "=MAX(B87:B91)"
$ws.Range("F87").Formula = "=AVERAGE(B87:B91)"
$ws.Range("G87").Formula = "=MIN(C87:C91)"
$ws.Range("H87").Formula = "=MAX(C87:C91)"
$ws.Range("I87").Formula = "=AVERAGE(C87:C91)"

# Row 88: Run#2
$ws.Range("A88").Value = "Run#2"
$ws.Range("B88").Value = 74.44
$ws.Range("C88").Value = 74.78
$ws.Range("F88").Formula = "=MAX(F87-D87,E87-F87)"
$ws.Range("I88").Formula = "=MAX(H87-I87,I87-G87)"

# Row 89: Run#3
$ws.Range("A89").Value = "Run#3"
$ws.Range("B89").Value = 74.79
$ws.Range("C89").Value = 73

# Row 91: header row for "Contrastive VILT new way 0.01"
$ws.Range("A91").Value = "Type - Contrastive VILT new way 0.01"
$ws.Range("B91").Value = "Test acc."
$ws.Range("C91").Value = "Dev acc."
$ws.Range("D91").Value = "Min test"
$ws.Range("E91").Value = "Max test"
$ws.Range("F91").Value = "Mean test"
$ws.Range("G91").Value = "Min dev"
$ws.Range("H91").Value = "Max dev"
$ws.Range("I91").Value = "Mean dev"

# Row 92: Run#1
$ws.Range("A92").Value = "Run#1"
$ws.Range("B92").Value = 73.8
$ws.Range("C92").Value = 72.3
$ws.Range("D92").Formula = "=MIN(B92:B96)"
$ws.Range("E92").Formula = "=MAX(B92:B96)"
$ws.Range("F92").Formula = "=AVERAGE(B92:B96)"
$ws.Range("G92").Formula = "=MIN(C92:C96)"
$ws.Range("H92").Formula = "=MAX(C92:C96)"
$ws.Range("I92").Formula = "=AVERAGE(C92:C96)"

# Row 93: Run#2
$ws.Range("A93").Value = "Run#2"
$ws.Range("B93").Value = 75.24
$ws.Range("C93").Value = 72.9
$ws.Range("F93").Formula = "=MAX(F92-D92,E92-F92)"
$ws.Range("I93").Formula = "=MAX(H92-I92,I92-G92)"

# Row 94: Run#3
$ws.Range("A94").Value = "Run#3"
$ws.Range("B94").Value = 73.95
$ws.Range("C94").Value = 73
$ws.Range("H94").Formula = "=(F92+I92)/2"

# --- view state: active tab scrolled/selected near the newly-added data ---
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 79 } catch { }
$ws.Range("D92").Select()
